$wb = $excel.ActiveWorkbook
$ws8 = $wb.Worksheets.Item(8)
$ws8.Copy($null, $ws8)
$ws9 = $wb.Worksheets.Item(9)
$ws9.Name = "Nädal 9"
Write-Host $ws9.Name
